{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body only contains the \"_GoBack\" bookmark\n// (no visible text). Insert the signature text at the very start of that\n// paragraph, ahead of the bookmark, then add a new empty paragraph after it.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"Kristen J. Rials (11/12/19)\", \"Start\");\nlastParagraph.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark sits by itself in the last paragraph of the\n# document (the line under \"please type your name below:\"). Insert the\n# signature text immediately before that bookmark so the bookmark keeps\n# wrapping the insertion point, then start a new empty paragraph after it.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Range.InsertBefore(\"Kristen J. Rials (11/12/19)\")\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null\n"}
